# "Örnek 9 - Eğer - 2.xlsx" — add the IF() formulas for the SICAK/NORMAL
# table (C:E, rows 2-8) and fill in the student-info block (H2:H4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: stand-alone (non shared) formulas ------------------------------
$ws.Range("C2").Formula = '=IF(B2>24,"SICAK", "NORMAL")'
$ws.Range("D2").Formula = '=IF(B2=25,"BUGÜN HAVA ÇOK GÜZEL ILIK", "SICAK YA DA SOĞUK OLABİLİR")'
$ws.Range("E2").Formula = '=IF(B2<24,"HAVAYA DİKKAT EDİN","BUGÜN HAVA 25 DERECE")'

# --- Rows 3-8: filled down, so these become shared formulas ----------------
$ws.Range("C3:C8").Formula = '=IF(B3>24,"SICAK", "NORMAL")'
$ws.Range("D3:D8").Formula = '=IF(B3=25,"BUGÜN HAVA ÇOK GÜZEL ILIK", "SICAK YA DA SOĞUK OLABİLİR")'

# Column E was filled down in two passes: first rows 3-5, then rows 6-8 were
# re-typed with a (typo'd) "DİKKAY" formula.
$ws.Range("E3:E5").Formula = '=IF(B3<24,"HAVAYA DİKKAT EDİN","BUGÜN HAVA 25 DERECE")'
$ws.Range("E6:E8").Formula = '=IF(B6<24,"HAVAYA DİKKAY EDİN","BUGÜN HAVA 25 DERECE")'

# --- Student info block ------------------------------------------------------
$ws.Range("H2").Value = 20215070019
$ws.Range("H3").Value = "KÜBRA ÇABUK"
$ws.Range("H4").Value = "YBS"

# --- Column widths (D, E widened; new narrow spacer column F) --------------
$ws.Columns.Item(4).ColumnWidth = 30.022135416666668
$ws.Columns.Item(5).ColumnWidth = 24.736979166666668
$ws.Columns.Item(6).ColumnWidth = 1.4518229166666665

# --- Selection moved to H7 --------------------------------------------------
$ws.Range("H7").Select()
